$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74..195 down to 75..196.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new data point.
$ws.Range("A74").Value = 4
$ws.Range("B74").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C74").Value = "Los Lagos"
$ws.Range("D74").Value = 44557
$ws.Range("E74").Value = 10
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100102
$ws.Range("H74").Value = "Cítricos"
$ws.Range("I74").Value = 100102006
$ws.Range("J74").Value = "Pomelo"
$ws.Range("K74").Value = "Start Ruby"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 120
$ws.Range("N74").Value = 10000
$ws.Range("O74").Value = 11000
$ws.Range("P74").Value = 10500
$ws.Range("Q74").Value = "$/caja 14 kilos empedrada"
$ws.Range("R74").Value = "Región de O'Higgins"
$ws.Range("S74").Value = 750
$ws.Range("T74").Value = 14
